$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The exported report dropped the redundant "id" column (old column A) and
# shifted everything else one column to the left; "pay_list" (old column L)
# is now the last column, K. Dropping the old last column first means the
# header/value writes below land on the final A:K layout directly.
$ws.Columns.Item(12).Delete()

# --- Header row (row 1) ------------------------------------------------
$headers = @(
    "serial number",
    "activated date",
    "pay_day",
    "pre-pay day",
    "on pause",
    "phone",
    "telegram username",
    "username",
    "email",
    "password",
    "pay_list"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- New data rows 2-5 ---------------------------------------------------
# Every field is free-form text in the source export (dates as plain
# "MM/DD/YYYY" strings, codes that look numeric, phone numbers with a
# leading "+", etc.), so force text formatting before assigning to stop
# Excel from auto-coercing them into dates/numbers.
$data = @(
    @("KIT400848944", "11/18/2024", "11/19/2024", "11/20/2024", "-", "+7 959 502 9416", "@M98_FaZa", "Фаза", "84@zov.icu", "Password1236", "341525"),
    @("KIT400864544", "11/18/2024", "11/19/2024", "11/20/2024", "-", "+7 904 902 4696", "@tankist613", "Alex210949 Alex210949", "84@zov.icu", "Password1237", "341526"),
    @("KIT400594583", "11/18/2024", "11/19/2024", "11/20/2024", "-", "+nan", "nan", "Илья .", "84@zov.icu", "Password1238", "341527"),
    @("KIT400594582131234", "11/18/2024", "11/19/2024", "11/20/2024", "-", "+32142134213", "asdsa", "Илья .", "84@zov.icu", "Password1238", "341527")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    $rowRange = $ws.Range("A" + $excelRow + ":K" + $excelRow)
    $rowRange.NumberFormat = "@"
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $row[$c]
    }
}

Write-Host "Done. UsedRange: $($ws.UsedRange.Address())"
